# usbwledc-top-pos.xlsx edit:
#   - "fix resistor placement instead of caps, replace mic, fix count of components"
#   The microphone component (row 14, designator MK1) had the wrong part
#   (GSA4737 / Microphone-6pin) placed at the wrong position/rotation.
#   Replace it with the correct part (ICS-43434 / InvenSense_ICS-43434-6_3.5x2.65mm)
#   at the corrected Mid X / Mid Y / Rotation, matching the non-default
#   (Calibri) font style used for the corrected Val/Package cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 = MK1 (microphone). Correct the part Val/Package text.
$ws.Range("B14").Value = "ICS-43434"
$ws.Range("C14").Value = "InvenSense_ICS-43434-6_3.5x2.65mm"

# The corrected cells pick up an explicit Calibri font (new style),
# distinguishing them from the rest of the untouched rows.
$ws.Range("B14:C14").Font.Name = "Calibri"

# Correct the placement (Mid X, Mid Y, Rotation) for MK1.
$ws.Range("D14").Value = 32.6
$ws.Range("E14").Value = -38.3
$ws.Range("F14").Value = 180

# Update the sheet's last active selection, left by the author after editing.
[void]$ws.Range("F15").Select()
